$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (E, F) before the old "ExcludedStudies_Excel_File_names" column ---
# Old layout: A Name | B population_name | C pop_radio_button | D Study_Types |
#             E ExcludedStudies_Excel_File_names | F ExcludedStudies_Excel_Files |
#             G Override_ExcludedStudies_Excel_File_names | H Override_ExcludedStudies_Excel_Files |
#             I ExpectedFilenames
# New layout adds "slrtype" / "slrtype_Radio_button" columns after Study_Types, pushing the
# rest two columns to the right (G..K).
$ws.Range("E1:F1").EntireColumn.Insert()

# --- Rename the population headers ---
$ws.Range("B1").Value = "Population"
$ws.Range("C1").Value = "Population_Radio_button"

# --- New slrtype columns header ---
$ws.Range("E1").Value = "slrtype"
$ws.Range("F1").Value = "slrtype_Radio_button"

# --- New slrtype data rows (mirrors Study_Types column D, plus a derived radio-button label) ---
$ws.Range("E2").Value = "Clinical"
$ws.Range("F2").Value = "Clinical_radio_button"

$ws.Range("E3").Value = "Economic"
$ws.Range("F3").Value = "Economic_radio_button"

$ws.Range("E4").Value = "Quality of Life"
$ws.Range("F4").Value = "Quality of Life_radio_button"

$ws.Range("E5").Value = "Real-world Evidence"
$ws.Range("F5").Value = "Real-world Evidence_radio_button"

# --- Column widths: keep E/F matching the look of column D ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# --- View / selection updates ---
$ws.Range("E3:F5").Select()
